# Db load from excel optimization
# The quiz-results export contained a duplicate/stale row (ID 111, "deepas")
# that should not have been included in the DB load. Remove that row
# entirely so the remaining rows shift up and the sheet's used range
# shrinks by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("5:5").Delete()
